# Auto-generated script applying FFXIV leve-profit data refresh across sheets
# Each block updates currentAveragePrice / LevePrice / LeveProfit columns (H-N)
# for specific rows on specific sheets, matching the scheduled market-data refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 20
$ws.Cells.Item(20, 8).Value = 800
$ws.Cells.Item(20, 9).Value = 800
$ws.Cells.Item(20, 11).Value = 800
$ws.Cells.Item(20, 13).Value = -570
# ALC row 32
$ws.Cells.Item(32, 8).Value = 573.3333
$ws.Cells.Item(32, 9).Value = 500
$ws.Cells.Item(32, 11).Value = 500
$ws.Cells.Item(32, 13).Value = -174
# ALC row 35
$ws.Cells.Item(35, 8).Value = 800
$ws.Cells.Item(35, 9).Value = 800
$ws.Cells.Item(35, 11).Value = 800
$ws.Cells.Item(35, 13).Value = -421
# ALC row 137
$ws.Cells.Item(137, 8).Value = 14287623
$ws.Cells.Item(137, 9).Value = 29412618
$ws.Cells.Item(137, 10).Value = 2905.3333
$ws.Cells.Item(137, 11).Value = 88237854
$ws.Cells.Item(137, 12).Value = 8715.999899999999
$ws.Cells.Item(137, 13).Value = -88235304
$ws.Cells.Item(137, 14).Value = -13815.9999

$ws = $wb.Worksheets.Item("ARM")
# ARM row 32
$ws.Cells.Item(32, 8).Value = 5681.92
$ws.Cells.Item(32, 9).Value = 5340.75
$ws.Cells.Item(32, 10).Value = 13870
$ws.Cells.Item(32, 11).Value = 5340.75
$ws.Cells.Item(32, 12).Value = 13870
$ws.Cells.Item(32, 13).Value = -5053.75
$ws.Cells.Item(32, 14).Value = -14444
# ARM row 61
$ws.Cells.Item(61, 8).Value = 1792.4706
$ws.Cells.Item(61, 9).Value = 1522.4
$ws.Cells.Item(61, 10).Value = 2383.25
$ws.Cells.Item(61, 11).Value = 1522.4
$ws.Cells.Item(61, 12).Value = 2383.25
$ws.Cells.Item(61, 13).Value = -1310.4
$ws.Cells.Item(61, 14).Value = -2807.25
# ARM row 132
$ws.Cells.Item(132, 8).Value = 1509.1887
$ws.Cells.Item(132, 9).Value = 1027.0233
$ws.Cells.Item(132, 10).Value = 3582.5
$ws.Cells.Item(132, 11).Value = 3081.0699
$ws.Cells.Item(132, 12).Value = 10747.5
$ws.Cells.Item(132, 13).Value = -551.0699000000004
$ws.Cells.Item(132, 14).Value = -15807.5
# ARM row 136
$ws.Cells.Item(136, 8).Value = 1792.4706
$ws.Cells.Item(136, 9).Value = 1522.4
$ws.Cells.Item(136, 10).Value = 2383.25
$ws.Cells.Item(136, 11).Value = 4567.200000000001
$ws.Cells.Item(136, 12).Value = 7149.75
$ws.Cells.Item(136, 13).Value = -2017.200000000001
$ws.Cells.Item(136, 14).Value = -12249.75

$ws = $wb.Worksheets.Item("BSM")
# BSM row 141
$ws.Cells.Item(141, 8).Value = 46348.75
$ws.Cells.Item(141, 10).Value = 46348.75
$ws.Cells.Item(141, 12).Value = 46348.75
$ws.Cells.Item(141, 14).Value = -56708.75

$ws = $wb.Worksheets.Item("CRP")
# CRP row 31
$ws.Cells.Item(31, 8).Value = 17244444
$ws.Cells.Item(31, 9).Value = 33334438
$ws.Cells.Item(31, 10).Value = 5163.9287
$ws.Cells.Item(31, 11).Value = 33334438
$ws.Cells.Item(31, 12).Value = 5163.9287
$ws.Cells.Item(31, 13).Value = -33334143
$ws.Cells.Item(31, 14).Value = -5753.9287
# CRP row 34
$ws.Cells.Item(34, 8).Value = 17244444
$ws.Cells.Item(34, 9).Value = 33334438
$ws.Cells.Item(34, 10).Value = 5163.9287
$ws.Cells.Item(34, 11).Value = 33334438
$ws.Cells.Item(34, 12).Value = 5163.9287
$ws.Cells.Item(34, 13).Value = -33334236
$ws.Cells.Item(34, 14).Value = -5567.9287
# CRP row 58
$ws.Cells.Item(58, 8).Value = 1587.0746
$ws.Cells.Item(58, 9).Value = 1081.9
$ws.Cells.Item(58, 10).Value = 2335.4814
$ws.Cells.Item(58, 11).Value = 1081.9
$ws.Cells.Item(58, 12).Value = 2335.4814
$ws.Cells.Item(58, 13).Value = -878.9000000000001
$ws.Cells.Item(58, 14).Value = -2741.4814
# CRP row 86
$ws.Cells.Item(86, 8).Value = 18909.7
$ws.Cells.Item(86, 9).Value = 37392.43
$ws.Cells.Item(86, 10).Value = 2737.3125
$ws.Cells.Item(86, 11).Value = 37392.43
$ws.Cells.Item(86, 12).Value = 2737.3125
$ws.Cells.Item(86, 13).Value = -36269.43
$ws.Cells.Item(86, 14).Value = -4983.3125
# CRP row 89
$ws.Cells.Item(89, 8).Value = 18909.7
$ws.Cells.Item(89, 9).Value = 37392.43
$ws.Cells.Item(89, 10).Value = 2737.3125
$ws.Cells.Item(89, 11).Value = 186962.15
$ws.Cells.Item(89, 12).Value = 13686.5625
$ws.Cells.Item(89, 13).Value = -181346.15
$ws.Cells.Item(89, 14).Value = -24918.5625
# CRP row 132
$ws.Cells.Item(132, 8).Value = 2014.2703
$ws.Cells.Item(132, 9).Value = 1450.0952
$ws.Cells.Item(132, 10).Value = 2754.75
$ws.Cells.Item(132, 11).Value = 4350.2856
$ws.Cells.Item(132, 12).Value = 8264.25
$ws.Cells.Item(132, 13).Value = -1820.2856
$ws.Cells.Item(132, 14).Value = -13324.25
# CRP row 136
$ws.Cells.Item(136, 8).Value = 1587.0746
$ws.Cells.Item(136, 9).Value = 1081.9
$ws.Cells.Item(136, 10).Value = 2335.4814
$ws.Cells.Item(136, 11).Value = 3245.7
$ws.Cells.Item(136, 12).Value = 7006.4442
$ws.Cells.Item(136, 13).Value = -695.7000000000003
$ws.Cells.Item(136, 14).Value = -12106.4442

$ws = $wb.Worksheets.Item("CUL")
# CUL row 113
$ws.Cells.Item(113, 8).Value = 546
$ws.Cells.Item(113, 9).Value = 549
$ws.Cells.Item(113, 10).Value = 545.4545000000001
$ws.Cells.Item(113, 11).Value = 1647
$ws.Cells.Item(113, 12).Value = 1636.3635
$ws.Cells.Item(113, 13).Value = 523
$ws.Cells.Item(113, 14).Value = -5976.3635
# CUL row 129
$ws.Cells.Item(129, 8).Value = 861.44446
$ws.Cells.Item(129, 10).Value = 1516.5
$ws.Cells.Item(129, 12).Value = 4549.5
$ws.Cells.Item(129, 14).Value = -14549.5
# CUL row 134
$ws.Cells.Item(134, 8).Value = 3932.356
$ws.Cells.Item(134, 9).Value = 1511.6111
$ws.Cells.Item(134, 10).Value = 4995.122
$ws.Cells.Item(134, 11).Value = 4534.8333
$ws.Cells.Item(134, 12).Value = 14985.366
$ws.Cells.Item(134, 13).Value = 535.1666999999998
$ws.Cells.Item(134, 14).Value = -25125.366
# CUL row 139
$ws.Cells.Item(139, 8).Value = 2904.1292
$ws.Cells.Item(139, 9).Value = 1696.5714
$ws.Cells.Item(139, 10).Value = 5440
$ws.Cells.Item(139, 11).Value = 5089.7142
$ws.Cells.Item(139, 12).Value = 16320
$ws.Cells.Item(139, 13).Value = 50.28579999999965
$ws.Cells.Item(139, 14).Value = -26600
# CUL row 140
$ws.Cells.Item(140, 8).Value = 1395.129
$ws.Cells.Item(140, 9).Value = 683.03845
$ws.Cells.Item(140, 10).Value = 5098
$ws.Cells.Item(140, 11).Value = 2049.11535
$ws.Cells.Item(140, 12).Value = 15294
$ws.Cells.Item(140, 13).Value = 3130.88465
$ws.Cells.Item(140, 14).Value = -25654

$ws = $wb.Worksheets.Item("GSM")
# GSM row 80
$ws.Cells.Item(80, 8).Value = 3439.6667
$ws.Cells.Item(80, 9).Value = 2225
$ws.Cells.Item(80, 10).Value = 3881.3635
$ws.Cells.Item(80, 11).Value = 2225
$ws.Cells.Item(80, 12).Value = 3881.3635
$ws.Cells.Item(80, 13).Value = -1227
$ws.Cells.Item(80, 14).Value = -5877.363499999999
# GSM row 83
$ws.Cells.Item(83, 8).Value = 3439.6667
$ws.Cells.Item(83, 9).Value = 2225
$ws.Cells.Item(83, 10).Value = 3881.3635
$ws.Cells.Item(83, 11).Value = 11125
$ws.Cells.Item(83, 12).Value = 19406.8175
$ws.Cells.Item(83, 13).Value = -6133
$ws.Cells.Item(83, 14).Value = -29390.8175
# GSM row 123
$ws.Cells.Item(123, 8).Value = 17322
$ws.Cells.Item(123, 10).Value = 17322
$ws.Cells.Item(123, 12).Value = 17322
$ws.Cells.Item(123, 14).Value = -22222

$ws = $wb.Worksheets.Item("LTW")
# LTW row 22
$ws.Cells.Item(22, 8).Value = 873.25
$ws.Cells.Item(22, 10).Value = 885.1429000000001
$ws.Cells.Item(22, 12).Value = 885.1429000000001
$ws.Cells.Item(22, 14).Value = -1475.1429
# LTW row 27
$ws.Cells.Item(27, 8).Value = 873.25
$ws.Cells.Item(27, 10).Value = 885.1429000000001
$ws.Cells.Item(27, 12).Value = 885.1429000000001
$ws.Cells.Item(27, 14).Value = -1099.1429
# LTW row 87
$ws.Cells.Item(87, 8).Value = 10000
$ws.Cells.Item(87, 9).Value = 10000
$ws.Cells.Item(87, 10).Value = 0
$ws.Cells.Item(87, 11).Value = 10000
$ws.Cells.Item(87, 12).ClearContents()
$ws.Cells.Item(87, 14).ClearContents()
$ws.Cells.Item(87, 13).Value = -8877
# LTW row 90
$ws.Cells.Item(90, 8).Value = 10000
$ws.Cells.Item(90, 9).Value = 10000
$ws.Cells.Item(90, 10).Value = 0
$ws.Cells.Item(90, 11).Value = 30000
$ws.Cells.Item(90, 12).ClearContents()
$ws.Cells.Item(90, 14).ClearContents()
$ws.Cells.Item(90, 13).Value = -24384

$ws = $wb.Worksheets.Item("WVR")
# WVR row 64
$ws.Cells.Item(64, 8).Value = 24347.883
$ws.Cells.Item(64, 10).Value = 24347.883
$ws.Cells.Item(64, 12).Value = 24347.883
$ws.Cells.Item(64, 14).Value = -24843.883
# WVR row 67
$ws.Cells.Item(67, 8).Value = 24347.883
$ws.Cells.Item(67, 10).Value = 24347.883
$ws.Cells.Item(67, 12).Value = 24347.883
$ws.Cells.Item(67, 14).Value = -26063.883
